$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("quiz")

# Update marking value (B11): 3 -> 5
$ws.Range("B11").Value = 5

# Update total correct marks (B12): 72 -> 120
$ws.Range("B12").Value = 120

# Update correct/total marks text (E12): 69/84 -> 120/140
$ws.Range("E12").Value = "120/140"
